$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D36").Value = "Unsupervised Semantic Segmentation"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/357"

$ws.Range("D51").Value = "[python+plotly] px.scatter 플롯에 추세선 넣기"
$ws.Range("E51").Value = "https://bskyvision.com/1255"
